$d = $word.ActiveDocument

# 1. Date change in the title paragraph: 05.06.24 -> 04.06.24
$d.Content.Find.Execute(
    "05.06.24", $false, $false, $false, $false, $false,
    $true, 1, $false, "04.06.24", 2) | Out-Null

# 2. Paper title paragraph: GraphAny... -> Are Emergent Abilities...
$d.Content.Find.Execute(
    "GraphAny: A Foundation Model for Node Classification on Any Graph",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "Are Emergent Abilities of Large Language Models a Mirage?", 2) | Out-Null

# 3. Third paragraph: replace whole paragraph body text
$p3 = $d.Paragraphs(3).Range
$p3.End = $p3.End - 1
$p3.Text = "היום המאמר שנסקור הוא מלפני שנה בערך והוא משך את תשומת ליבי בגלל שהוא חוקר מה שנקרא emergent capabilities של מודלי שפה - כלומר יכולתם ללמוד משימות חדשות. המאמר בוחן האם למודלי שפה אכן יש יכולת ללמוד משימות שהם אומנו עליהם בצורה מפורשת (פחות או יותר) או שזו אשליה הנובעת מאיך שאנו מודדים את היכולות האלו."

# 4. Fourth paragraph: replace whole paragraph body text
$p4 = $d.Paragraphs(4).Range
$p4.End = $p4.End - 1
$p4.Text = "מאמר: https://arxiv.org/abs/2304.15004"

# 5. Fifth paragraph: collapse the multi-run/break paragraph into a single
#    line with the updated telegram link only.
$p5 = $d.Paragraphs(5).Range
$p5.End = $p5.End - 1
$p5.Text = "טלגרם: https://t.me/MathyAIwithMike/76"
